# Updates cryptos list figures (price/volume) and reorders a few coin rows,
# matching the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" on numeric-looking strings forces Excel to keep them as text
# (matching the original inlineStr cell type) instead of auto-converting to a number.

$ws.Cells.Item(2, 4).Value = "30.633.74"
$ws.Cells.Item(2, 5).Value = "  +0.49%  "
$ws.Cells.Item(3, 4).Value = "1.964.29"
$ws.Cells.Item(3, 5).Value = "  +2.67%  "
$ws.Cells.Item(4, 4).Value = "'0.9989"
$ws.Cells.Item(4, 5).Value = "  +0.01%  "
$ws.Cells.Item(5, 4).Value = "'249.00"
$ws.Cells.Item(5, 5).Value = "  +1.49%  "
$ws.Cells.Item(6, 4).Value = "'0.9987"
$ws.Cells.Item(6, 5).Value = "  -0.03%  "
$ws.Cells.Item(7, 4).Value = "'0.4802"
$ws.Cells.Item(7, 5).Value = "  -0.48%  "
$ws.Cells.Item(8, 4).Value = "'0.2943"
$ws.Cells.Item(8, 5).Value = "  +1.62%  "
$ws.Cells.Item(9, 4).Value = "'0.06799"
$ws.Cells.Item(9, 5).Value = "  +1.11%  "
$ws.Cells.Item(10, 4).Value = "'111.99"
$ws.Cells.Item(10, 5).Value = "  +1.04%  "
$ws.Cells.Item(11, 4).Value = "'19.39"
$ws.Cells.Item(11, 5).Value = "  +0.96%  "
$ws.Cells.Item(12, 4).Value = "1.931.54"
$ws.Cells.Item(12, 5).Value = "  +0.95%  "
$ws.Cells.Item(13, 4).Value = "'0.07696"
$ws.Cells.Item(13, 5).Value = "  +1.85%  "
$ws.Cells.Item(14, 4).Value = "'5.487"
$ws.Cells.Item(14, 5).Value = "  +3.99%  "
$ws.Cells.Item(15, 4).Value = "'0.6876"
$ws.Cells.Item(15, 5).Value = "  +2.39%  "
$ws.Cells.Item(16, 4).Value = "'294.90"
$ws.Cells.Item(16, 5).Value = "  +1.97%  "
$ws.Cells.Item(17, 4).Value = "30.651.94"
$ws.Cells.Item(17, 5).Value = "  +0.58%  "
$ws.Cells.Item(18, 4).Value = "'13.29"
$ws.Cells.Item(18, 5).Value = "  +3.19%  "
$ws.Cells.Item(19, 4).Value = "'5.662"
$ws.Cells.Item(19, 5).Value = "  +3.49%  "
$ws.Cells.Item(20, 4).Value = "2.219.38"
$ws.Cells.Item(20, 5).Value = "  +2.53%  "
$ws.Cells.Item(21, 4).Value = "'0.000007666"
$ws.Cells.Item(21, 5).Value = "  +0.88%  "
$ws.Cells.Item(22, 4).Value = "'0.9981"
$ws.Cells.Item(22, 5).Value = "  -0.08%  "
$ws.Cells.Item(23, 4).Value = "'1.001"
$ws.Cells.Item(23, 5).Value = "  +0.21%  "
$ws.Cells.Item(24, 4).Value = "'6.602"
$ws.Cells.Item(24, 5).Value = "  +3.12%  "
$ws.Cells.Item(25, 4).Value = "'9.740"
$ws.Cells.Item(25, 5).Value = "  +2.88%  "
$ws.Cells.Item(26, 4).Value = "'168.79"
$ws.Cells.Item(26, 5).Value = "  +2.86%  "
$ws.Cells.Item(27, 4).Value = "'20.31"
$ws.Cells.Item(27, 5).Value = "  -0.28%  "
$ws.Cells.Item(28, 4).Value = "'2.207"
$ws.Cells.Item(28, 5).Value = "  +3.53%  "
$ws.Cells.Item(29, 5).Value = "  +3.16%  "
$ws.Cells.Item(30, 4).Value = "'1.433"
$ws.Cells.Item(30, 5).Value = "  +1.91%  "
$ws.Cells.Item(31, 4).Value = "'4.694"
$ws.Cells.Item(31, 5).Value = "  +15.63%  "
$ws.Cells.Item(32, 4).Value = "'4.419"
$ws.Cells.Item(32, 5).Value = "  +5.68%  "
$ws.Cells.Item(33, 4).Value = "'0.05079"
$ws.Cells.Item(33, 5).Value = "  +1.72%  "
$ws.Cells.Item(34, 4).Value = "'0.7775"
$ws.Cells.Item(34, 5).Value = "  +6.55%  "
$ws.Cells.Item(35, 4).Value = "'1.168"
$ws.Cells.Item(35, 5).Value = "  +2.96%  "
$ws.Cells.Item(36, 4).Value = "'0.02066"
$ws.Cells.Item(36, 5).Value = "  +0.78%  "
$ws.Cells.Item(37, 4).Value = "'2.727"
$ws.Cells.Item(37, 5).Value = "  +0.33%  "
$ws.Cells.Item(38, 4).Value = "'2.715"
$ws.Cells.Item(38, 5).Value = "  +1.85%  "
$ws.Cells.Item(39, 4).Value = "'2.061"
$ws.Cells.Item(39, 5).Value = "  +2.36%  "
$ws.Cells.Item(40, 4).Value = "'110.99"
$ws.Cells.Item(40, 5).Value = "  +0.23%  "
$ws.Cells.Item(41, 4).Value = "'0.4463"
$ws.Cells.Item(41, 5).Value = "  +0.49%  "
$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(42, 4).Value = "'6.039"
$ws.Cells.Item(42, 5).Value = "  +3.36%  "
$ws.Cells.Item(43, 2).Value = "TrustWalletToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(43, 4).Value = "'0.8736"
$ws.Cells.Item(43, 5).Value = "  +0.65%  "
$ws.Cells.Item(44, 2).Value = "PaxDollar"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(44, 4).Value = "'0.9995"
$ws.Cells.Item(44, 5).Value = "  +0.07%  "
$ws.Cells.Item(45, 2).Value = "Aave"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(45, 4).Value = "'69.61"
$ws.Cells.Item(45, 5).Value = "  +2.05%  "
$ws.Cells.Item(46, 4).Value = "'7.394"
$ws.Cells.Item(46, 5).Value = "  +0.37%  "
$ws.Cells.Item(47, 4).Value = "'9.329"
$ws.Cells.Item(47, 5).Value = "  +0.23%  "
$ws.Cells.Item(48, 4).Value = "'0.1253"
$ws.Cells.Item(48, 5).Value = "  +1.26%  "
$ws.Cells.Item(49, 4).Value = "'47.99"
$ws.Cells.Item(49, 5).Value = "  -2.16%  "
$ws.Cells.Item(50, 4).Value = "'35.65"
$ws.Cells.Item(50, 5).Value = "  +2.26%  "
$ws.Cells.Item(51, 2).Value = "NEARProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(51, 4).Value = "'1.470"
$ws.Cells.Item(51, 5).Value = "  +1.57%  "
